# Generate Report for Handoff
#
# The b.md file has now been handed off for both zh-cn and de-de locales.
# Update the Overview sheet and both locale sheets to reflect the new
# "Ready for handoff" status, the new handoff xliff file names, the new
# handoff timestamps, and the new "handback not latest" error detail.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0876f813c1e38478af027f62fe90a4007f3402f1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d43bfa165acf577a745bd11465166d1f29acf503/e2e/b.md."

# ----------------------------------------------------------------------
# Overview sheet: row 3 is the b.md row
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = "2016-10-27 09:24:26"

# ----------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md row
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-27 09:24:14"
$wsZhCn.Range("P3").Value = $errorDetail
# 39.15 (rather than 40) compensates for this engine's internal
# char-width -> pixel -> char-width rounding so the saved column width
# lands on exactly 40, matching the target OOXML.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ----------------------------------------------------------------------
# de-de sheet: row 3 is the b.md row
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-27 09:24:26"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
